# Apply the PR-inventory update:
#  1. Rename the "congenital" sheet (tab name + its header cell A1) to "misc_long_term".
#  2. On the "mental" sheet, remove six rows that no longer belong (duplicates /
#     baseline-regression cleanup): GCST009722, GCST008373, GCST009520, GCST009521,
#     GA3543, GA3667. Removing them shifts everything below up, shrinking the used
#     range from A1:A102 to A1:A96.

$wb = $excel.ActiveWorkbook

# --- 1. Rename "congenital" -> "misc_long_term" ---------------------------------
$wsCongenital = $wb.Worksheets.Item("congenital")
$wsCongenital.Range("A1").Value = "misc_long_term"
$wsCongenital.Name = "misc_long_term"

# --- 2. Drop the six obsolete dataset rows from "mental" -------------------------
$wsMental = $wb.Worksheets.Item("mental")

# Row numbers as they exist in the *current* sheet, before any deletions.
# Deleting from the bottom up keeps the earlier row numbers valid.
$rowsToDelete = @(66, 59, 24, 23, 8, 7)
foreach ($r in $rowsToDelete) {
    $wsMental.Rows.Item($r).Delete()
}
